{"js": "// The \"Note\" paragraph currently ends with a run of trailing spaces, and\n// the hidden \"_GoBack\" bookmark sits at the end of the following\n// paragraph (right after \"...do your own research.      \"). This edit:\n//   1. moves the \"_GoBack\" bookmark to the end of the \"Note\" paragraph\n//      (immediately after the word \"Note\"); and\n//   2. removes the run of trailing spaces that used to follow \"Note\".\n// All of the runs/text in the following \"This tutorial...\" paragraph are\n// left untouched other than losing the bookmark that used to trail them.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose text is \"Note\" plus optional trailing\n// whitespace (robust to the exact run/whitespace split).\nlet noteParagraph = null;\nfor (const p of paragraphs.items) {\n  if (/^Note\\s*$/.test(p.text)) {\n    noteParagraph = p;\n    break;\n  }\n}\nif (!noteParagraph) {\n  throw new Error('Could not locate the \"Note\" paragraph.');\n}\n\n// 1) Remove the old \"_GoBack\" bookmark wherever it currently is. This\n//    must happen BEFORE inserting the new one at its target location,\n//    because this host does not relocate a bookmark when a new one is\n//    inserted under the same (already-used) name.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Insert \"_GoBack\" at the end of the \"Note\" paragraph (still before\n//    the trailing-space run is removed, but that does not matter since\n//    we are anchoring on the paragraph's end, not on specific text).\nconst noteEnd = noteParagraph.getRange(\"End\");\nnoteEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Delete the trailing run of spaces that follows \"Note\" (whatever its\n//    length): locate the literal word \"Note\" and take everything from\n//    right after it through to the end of the paragraph.\nconst noteWordResults = noteParagraph.search(\"Note\", { matchWildcards: false });\nnoteWordResults.load(\"items\");\nawait context.sync();\nif (noteWordResults.items.length > 0) {\n  const afterNote = noteWordResults.items[0].getRange(\"End\");\n  const paragraphEnd = noteParagraph.getRange(\"End\");\n  const trailingSpaces = afterNote.expandTo(paragraphEnd);\n  trailingSpaces.load(\"text\");\n  await context.sync();\n  if (trailingSpaces.text.length > 0) {\n    trailingSpaces.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# The \"Note\" paragraph currently ends with a run of trailing spaces, and\n# the hidden \"_GoBack\" bookmark sits at the end of the following\n# paragraph (right after \"...do your own research.      \"). This edit:\n#   1. moves the \"_GoBack\" bookmark to the end of the \"Note\" paragraph\n#      (immediately after the word \"Note\"); and\n#   2. removes the run of trailing spaces that used to follow \"Note\".\n# All of the runs/text in the following \"This tutorial...\" paragraph are\n# left untouched other than losing the bookmark that used to trail them.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose text is \"Note\" plus trailing whitespace.\n$noteParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd() -eq \"Note\") {\n        $noteParagraph = $p\n        break\n    }\n}\nif ($noteParagraph -eq $null) {\n    throw \"Could not locate the 'Note' paragraph.\"\n}\n\n# 1) Move the \"_GoBack\" bookmark to the end of the \"Note\" paragraph\n#    BEFORE deleting the trailing spaces, so the anchor position is\n#    unambiguous (computing it only from positions left behind by a\n#    prior delete is unreliable in this host).\n$noteRange = $d.Content\n$noteRange.Find.Execute(\"Note\") | Out-Null\n$noteRange.Collapse(0)  # wdCollapseEnd -> right after \"Note\"\n$d.Bookmarks.Add(\"_GoBack\", $noteRange)\n\n# 2) Delete the trailing run of spaces that used to follow \"Note\".\n$spaceRange = $noteParagraph.Range.Duplicate\n$spaceRange.Find.Execute(\"     \") | Out-Null\n$spaceRange.Delete()\n"}
